$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 396.4
$ws.Cells.Item(43, 9).Value = 310.4
$ws.Cells.Item(43, 10).Value = 439.4
$ws.Cells.Item(43, 11).Value = 310.4
$ws.Cells.Item(43, 12).Value = 439.4
$ws.Cells.Item(43, 13).Value = -241.4
$ws.Cells.Item(43, 14).Value = -577.4

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 225463.52
$ws.Cells.Item(98, 9).Value = 256002.75
$ws.Cells.Item(98, 10).Value = 1509.1666
$ws.Cells.Item(98, 11).Value = 256002.75
$ws.Cells.Item(98, 12).Value = 1509.1666
$ws.Cells.Item(98, 13).Value = -254504.75
$ws.Cells.Item(98, 14).Value = -4505.1666

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 397297.5
$ws.Cells.Item(107, 9).Value = 617706.25
$ws.Cells.Item(107, 10).Value = 561.7
$ws.Cells.Item(107, 11).Value = 617706.25
$ws.Cells.Item(107, 12).Value = 561.7
$ws.Cells.Item(107, 13).Value = -615786.25
$ws.Cells.Item(107, 14).Value = -4401.7

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 225463.52
$ws.Cells.Item(122, 9).Value = 256002.75
$ws.Cells.Item(122, 10).Value = 1509.1666
$ws.Cells.Item(122, 11).Value = 768008.25
$ws.Cells.Item(122, 12).Value = 4527.4998
$ws.Cells.Item(122, 13).Value = -765558.25
$ws.Cells.Item(122, 14).Value = -9427.4998

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1068.5769
$ws.Cells.Item(129, 10).Value = 1287.2106
$ws.Cells.Item(129, 12).Value = 3861.6318
$ws.Cells.Item(129, 14).Value = -13861.6318

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15167.156
$ws.Cells.Item(32, 9).Value = 1326.5286
$ws.Cells.Item(32, 11).Value = 1326.5286
$ws.Cells.Item(32, 13).Value = -1039.5286

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1446.125
$ws.Cells.Item(45, 9).Value = 1129.4546
$ws.Cells.Item(45, 10).Value = 2142.8
$ws.Cells.Item(45, 11).Value = 1129.4546
$ws.Cells.Item(45, 12).Value = 2142.8
$ws.Cells.Item(45, 13).Value = -752.4546
$ws.Cells.Item(45, 14).Value = -2896.8

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 22229528
$ws.Cells.Item(97, 9).Value = 27786410
$ws.Cells.Item(97, 10).Value = 1999
$ws.Cells.Item(97, 11).Value = 27786410
$ws.Cells.Item(97, 12).Value = 1999
$ws.Cells.Item(97, 13).Value = -27785914
$ws.Cells.Item(97, 14).Value = -2991

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 33399.8
$ws.Cells.Item(133, 10).Value = 33399.8
$ws.Cells.Item(133, 12).Value = 33399.8
$ws.Cells.Item(133, 14).Value = -38459.8

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 811.3
$ws.Cells.Item(94, 9).Value = 811.3
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 811.3
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).Value = -360.3
$ws.Cells.Item(94, 13).ClearContents()

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 41525
$ws.Cells.Item(138, 10).Value = 41525
$ws.Cells.Item(138, 12).Value = 41525
$ws.Cells.Item(138, 14).Value = -51805

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 35500
$ws.Cells.Item(140, 10).Value = 35500
$ws.Cells.Item(140, 12).Value = 35500
$ws.Cells.Item(140, 14).Value = -45860

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1460.8636
$ws.Cells.Item(31, 9).Value = 1177.1177
$ws.Cells.Item(31, 10).Value = 2425.6
$ws.Cells.Item(31, 11).Value = 1177.1177
$ws.Cells.Item(31, 12).Value = 2425.6
$ws.Cells.Item(31, 13).Value = -882.1177
$ws.Cells.Item(31, 14).Value = -3015.6

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1460.8636
$ws.Cells.Item(34, 9).Value = 1177.1177
$ws.Cells.Item(34, 10).Value = 2425.6
$ws.Cells.Item(34, 11).Value = 1177.1177
$ws.Cells.Item(34, 12).Value = 2425.6
$ws.Cells.Item(34, 13).Value = -975.1177
$ws.Cells.Item(34, 14).Value = -2829.6

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 7814489
$ws.Cells.Item(99, 10).Value = 2760
$ws.Cells.Item(99, 12).Value = 2760
$ws.Cells.Item(99, 14).Value = -5756

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 458.625
$ws.Cells.Item(107, 10).Value = 999.6667
$ws.Cells.Item(107, 12).Value = 999.6667
$ws.Cells.Item(107, 14).Value = -4839.6667

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 7814489
$ws.Cells.Item(126, 10).Value = 2760
$ws.Cells.Item(126, 12).Value = 8280
$ws.Cells.Item(126, 14).Value = -13220

# CRP row 137
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(137, 8).Value = 21390
$ws.Cells.Item(137, 10).Value = 32780
$ws.Cells.Item(137, 12).Value = 32780
$ws.Cells.Item(137, 14).Value = -42980

# CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 47460
$ws.Cells.Item(138, 10).Value = 47460
$ws.Cells.Item(138, 12).Value = 47460
$ws.Cells.Item(138, 14).Value = -57740

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 20892
$ws.Cells.Item(3, 9).Value = 20892
$ws.Cells.Item(3, 11).Value = 62676
$ws.Cells.Item(3, 13).Value = -62564

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 56444.445
$ws.Cells.Item(37, 10).Value = 56444.445
$ws.Cells.Item(37, 12).Value = 169333.335
$ws.Cells.Item(37, 14).Value = -169557.335

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 6173813.5
$ws.Cells.Item(132, 9).Value = 739.3
$ws.Cells.Item(132, 10).Value = 9805034
$ws.Cells.Item(132, 11).Value = 6653.7
$ws.Cells.Item(132, 12).Value = 88245306
$ws.Cells.Item(132, 13).Value = -4123.7
$ws.Cells.Item(132, 14).Value = -88250366

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 5955
$ws.Cells.Item(133, 9).Value = 865
$ws.Cells.Item(133, 10).Value = 8500
$ws.Cells.Item(133, 11).Value = 2595
$ws.Cells.Item(133, 12).Value = 25500
$ws.Cells.Item(133, 13).Value = 2465
$ws.Cells.Item(133, 14).Value = -35620

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 5076.6
$ws.Cells.Item(134, 9).Value = 3447.5715
$ws.Cells.Item(134, 11).Value = 10342.7145
$ws.Cells.Item(134, 13).Value = -5272.7145

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 4216.5483
$ws.Cells.Item(136, 9).Value = 2152
$ws.Cells.Item(136, 10).Value = 4613.577
$ws.Cells.Item(136, 11).Value = 6456
$ws.Cells.Item(136, 12).Value = 13840.731
$ws.Cells.Item(136, 13).Value = -1356
$ws.Cells.Item(136, 14).Value = -24040.731

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 972.5
$ws.Cells.Item(138, 9).Value = 942.63635
$ws.Cells.Item(138, 10).Value = 1301
$ws.Cells.Item(138, 11).Value = 2827.90905
$ws.Cells.Item(138, 12).Value = 3903
$ws.Cells.Item(138, 13).Value = 2312.09095
$ws.Cells.Item(138, 14).Value = -14183

# GSM row 55
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 3210
$ws.Cells.Item(55, 10).Value = 6000
$ws.Cells.Item(55, 12).Value = 6000
$ws.Cells.Item(55, 14).Value = -6654

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1833.1666
$ws.Cells.Item(113, 9).Value = 1999.5
$ws.Cells.Item(113, 10).Value = 1750
$ws.Cells.Item(113, 11).Value = 1999.5
$ws.Cells.Item(113, 12).Value = 1750
$ws.Cells.Item(113, 13).Value = 170.5
$ws.Cells.Item(113, 14).Value = -6090

# GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 100011896
$ws.Cells.Item(135, 10).Value = 100011896
$ws.Cells.Item(135, 12).Value = 100011896
$ws.Cells.Item(135, 14).Value = -100022036

# GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(137, 8).Value = 43499.75
$ws.Cells.Item(137, 10).Value = 43499.75
$ws.Cells.Item(137, 12).Value = 43499.75
$ws.Cells.Item(137, 14).Value = -53699.75

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 42674.875
$ws.Cells.Item(138, 10).Value = 42674.875
$ws.Cells.Item(138, 12).Value = 42674.875
$ws.Cells.Item(138, 14).Value = -52954.875

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 384.3889
$ws.Cells.Item(55, 9).Value = 306.9091
$ws.Cells.Item(55, 10).Value = 506.14285
$ws.Cells.Item(55, 11).Value = 306.9091
$ws.Cells.Item(55, 12).Value = 506.14285
$ws.Cells.Item(55, 13).Value = -133.9091
$ws.Cells.Item(55, 14).Value = -852.14285

# LTW row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(125, 8).Value = 45500
$ws.Cells.Item(125, 10).Value = 45500
$ws.Cells.Item(125, 12).Value = 45500
$ws.Cells.Item(125, 14).Value = -55340

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 43928.75
$ws.Cells.Item(127, 10).Value = 43928.75
$ws.Cells.Item(127, 12).Value = 43928.75
$ws.Cells.Item(127, 14).Value = -53848.75

# LTW row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(128, 8).Value = 68333.336
$ws.Cells.Item(128, 10).Value = 68333.336
$ws.Cells.Item(128, 12).Value = 68333.336
$ws.Cells.Item(128, 14).Value = -78293.336

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 47941.11
$ws.Cells.Item(133, 10).Value = 47941.11
$ws.Cells.Item(133, 12).Value = 47941.11
$ws.Cells.Item(133, 14).Value = -53001.11

# LTW row 135
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(135, 8).Value = 32428
$ws.Cells.Item(135, 10).Value = 32428
$ws.Cells.Item(135, 12).Value = 32428
$ws.Cells.Item(135, 14).Value = -42568

# LTW row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(141, 8).Value = 45500
$ws.Cells.Item(141, 10).Value = 45500
$ws.Cells.Item(141, 12).Value = 45500
$ws.Cells.Item(141, 14).Value = -55860

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 125574.164
$ws.Cells.Item(135, 9).Value = 20000
$ws.Cells.Item(135, 10).Value = 146689
$ws.Cells.Item(135, 11).Value = 20000
$ws.Cells.Item(135, 12).Value = 146689
$ws.Cells.Item(135, 13).Value = -14930
$ws.Cells.Item(135, 14).Value = -156829

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 44400
$ws.Cells.Item(139, 10).Value = 44400
$ws.Cells.Item(139, 12).Value = 44400
$ws.Cells.Item(139, 14).Value = -54680

# WVR row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(141, 8).Value = 68857.5
$ws.Cells.Item(141, 10).Value = 68857.5
$ws.Cells.Item(141, 12).Value = 68857.5
$ws.Cells.Item(141, 14).Value = -79217.5
